# Revert responsive design implementation
#
# Appends 3 new sensor-data rows (29, 30, 31) to the bottom of both the
# "ROW35-FE-LIFTER" and "ROW35-MID-LIFTER" worksheets, matching the rows
# already present on "ROW02-FE-LIFTER" / "ROW02-MID-LIFTER". The worksheet
# dimension (A1:I28 -> A1:I31) updates automatically as the new cells are
# written.

$wb = $excel.ActiveWorkbook

# Column A datetime values are numbers formatted as "YYYY-MM-DD HH:MM:SS"
# (same custom number format already used by the existing rows).
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# G column holds a very large number that can't be written as a PowerShell
# numeric literal directly (overflows integer parsing before the [double]
# cast applies), so build it by casting the string form instead.
$gValue = [double]"5.686312626471138e+23"

# NOTE: this interpreter does not reliably bind PowerShell named parameters
# (`-Name value`) on custom functions, so this helper uses positional
# parameters only.
function Add-SensorRow($Worksheet, $RowIndex, $AValue, $BValue, $CValue, $DValue, $EValue, $FValue, $GValue, $HValue, $IValue) {
    $Worksheet.Range("A$RowIndex").Value = $AValue
    $Worksheet.Range("A$RowIndex").NumberFormat = $dateFormat
    $Worksheet.Range("B$RowIndex").Value = $BValue
    $Worksheet.Range("C$RowIndex").Value = $CValue
    $Worksheet.Range("D$RowIndex").Value = $DValue
    $Worksheet.Range("E$RowIndex").Value = $EValue
    $Worksheet.Range("F$RowIndex").Value = $FValue
    $Worksheet.Range("G$RowIndex").Value = $GValue
    $Worksheet.Range("H$RowIndex").Value = $HValue
    $Worksheet.Range("I$RowIndex").Value = $IValue
}

# ---------------------------------------------------------------------
# ROW35-FE-LIFTER
# ---------------------------------------------------------------------
$wsFe = $wb.Worksheets.Item("ROW35-FE-LIFTER")

$feRows = @(
    @{
        Row = 29
        A = 45729.73239443287
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xd"
        F = 400
        H = 400
        I = 13
    },
    @{
        Row = 30
        A = 45729.73241640046
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xd"
        F = 400
        H = 400
        I = 13
    },
    @{
        Row = 31
        A = 45729.73243972223
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xd"
        F = 400
        H = 400
        I = 13
    }
)

foreach ($r in $feRows) {
    Add-SensorRow $wsFe $r.Row $r.A $r.B $r.C $r.D $r.E $r.F $gValue $r.H $r.I
}

# ---------------------------------------------------------------------
# ROW35-MID-LIFTER
# ---------------------------------------------------------------------
$wsMid = $wb.Worksheets.Item("ROW35-MID-LIFTER")

$midRows = @(
    @{
        Row = 29
        A = 45729.58037369213
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x86,"
        E = "0x4"
        F = 400
        H = 390
        I = 4
    },
    @{
        Row = 30
        A = 45729.58039555555
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x86,"
        E = "0x4"
        F = 400
        H = 390
        I = 4
    },
    @{
        Row = 31
        A = 45729.58041870371
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x86,"
        E = "0x4"
        F = 400
        H = 390
        I = 4
    }
)

foreach ($r in $midRows) {
    Add-SensorRow $wsMid $r.Row $r.A $r.B $r.C $r.D $r.E $r.F $gValue $r.H $r.I
}
